$d = $word.ActiveDocument

# The document currently ends with a list paragraph containing
# "The cat cannot be left ... across the river. " followed by the
# hidden _GoBack bookmark. We need to add a brand-new list paragraph
# right after it with the "Possible solution: ..." text, and the
# _GoBack bookmark should end up tracking the end of that new text
# (mirroring how Word moves _GoBack to the last edited location).

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)

# Create the new paragraph (inherits list formatting from the paragraph
# it was split from, matching the target numbering/style).
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($lastParaIndex + 1)

$solutionText = "Possible solution: First take over the parrot, and then go get the cat, however when you drop the cat off, take the parrot back to the original side. Drop the parrot off at the original side and take the bag of seed over to the destination side and leave the seed with the cat. Then go back across and pick the parrot back up and return to the destination. "

# Append a one-character buffer so that, while we are positioning the
# bookmark, the insertion point is never the paragraph's very last
# interior offset (immediately before the paragraph mark) -- placing a
# zero-length bookmark exactly there mis-resolves in this runtime, so
# we dodge it by bookmarking just before a throw-away trailing
# character and then deleting that character afterwards.
$newPara.Range.Text = $solutionText + "#"

$newPara = $d.Paragraphs.Item($lastParaIndex + 1)
$rng = $newPara.Range
$bookmarkPos = $rng.End - 2

$point = $d.Range($bookmarkPos, $bookmarkPos)

# Move (recreate) the _GoBack bookmark at the end of the new text.
$d.Bookmarks.Item("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $point)

# Remove the throw-away buffer character now that the bookmark has
# been anchored just before it -- the bookmark (a collapsed point)
# stays put while the character after it disappears.
$bufferRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$bufferRange.Delete()
